# Adder - front panel PCB
# Extend the jack/LED spacing sheet with the absolute X positions of every
# jack/LED/mounting hole along the panel, plus the derived centers/widths
# of the cut-out rectangles ("rects") used for the front-panel artwork.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- remove the old scratch labels that used to sit in D2:E2 -------------
$ws.Range("D2:E2").ClearContents() | Out-Null

# --- column A/B : jack-to-jack spacing table (unchanged) ------------------
$ws.Range("A1").Value = "top jack"
$ws.Range("B1").Value = 1.65

$ws.Range("A2").Value = "jack-jack"
$ws.Range("B2").Value = 0.42500000000000027

$ws.Range("A3").Value = "jack-LED"
$ws.Range("B3").Value = 0.37000000000000011

$ws.Range("A4").Value = "LED-jack"
$ws.Range("B4").Value = 0.37999999999999989

$ws.Range("A5").Value = "J4-J5"
$ws.Range("B5").Value = 0.64999999999999991

$ws.Range("A6").Value = "bottom jack"
$ws.Range("B6").Value = 5.5

$ws.Range("A8").Value = "J1"
$ws.Range("B8").Formula = '=$B$1'

$ws.Range("A9").Value = "J2"
$ws.Range("B9").Formula = '=B8+$B$2'

$ws.Range("A10").Value = "J3"
$ws.Range("B10").Formula = '=B9+$B$2'

$ws.Range("A11").Value = "D1"
$ws.Range("B11").Formula = '=B10+$B$3'

$ws.Range("A12").Value = "J4"
$ws.Range("B12").Formula = '=B11+$B$4'

$ws.Range("A13").Value = "J5"
$ws.Range("B13").Formula = '=B12+$B$5'

$ws.Range("A14").Value = "J6"
$ws.Range("B14").Formula = '=B13+$B$2'

$ws.Range("A15").Value = "J7"
$ws.Range("B15").Formula = '=B14+$B$2'

$ws.Range("A16").Value = "D2"
$ws.Range("B16").Formula = '=B15+$B$3'

$ws.Range("A17").Value = "J8"
$ws.Range("B17").Formula = '=B16+$B$4'
$ws.Range("C17").Formula = '=B17=B6'

# --- columns D/E/F : absolute X position of every feature + its pitch ----
$ws.Range("D1").Value = "J1"
$ws.Range("E1").Value = 65.227999999999994

$ws.Range("F2").Formula = '=E3-E1'

$ws.Range("D3").Value = "D1"
$ws.Range("E3").Value = 74.626000000000005

$ws.Range("F4").Formula = '=E5-E3'

$ws.Range("D5").Value = "J2"
$ws.Range("E5").Value = 84.278000000000006

$ws.Range("F6").Formula = '=E7-E5'

$ws.Range("D7").Value = "J3"
$ws.Range("E7").Value = 95.072999999999993

$ws.Range("F8").Formula = '=E9-E7'

$ws.Range("D9").Value = "J4"
$ws.Range("E9").Value = 105.86799999999999

$ws.Range("F10").Formula = '=E11-E9'

$ws.Range("D11").Value = "J5"
$ws.Range("E11").Value = 122.378

$ws.Range("F12").Formula = '=E13-E11'

$ws.Range("D13").Value = "D2"
$ws.Range("E13").Value = 131.77600000000001

$ws.Range("F14").Formula = '=E15-E13'

$ws.Range("D15").Value = "J6"
$ws.Range("E15").Value = 141.428

$ws.Range("F16").Formula = '=E17-E15'

$ws.Range("D17").Value = "J7"
$ws.Range("E17").Value = 152.22300000000001

$ws.Range("F18").Formula = '=E19-E17'

$ws.Range("D19").Value = "J8"
$ws.Range("E19").Value = 163.018

# --- columns G/H/I : cumulative X position rebuilt from the pitches ------
$ws.Range("G1").Value = "J1"
$ws.Range("H1").Value = 65.227999999999994
$ws.Range("I2").Value = 10.794999999999987

$ws.Range("G3").Value = "J2"
$ws.Range("H3").Formula = '=H1+I2'
$ws.Range("I4").Value = 10.794999999999987

$ws.Range("G5").Value = "J3"
$ws.Range("H5").Formula = '=H3+I4'
$ws.Range("I6").Value = 9.3980000000000103

$ws.Range("G7").Value = "D1"
$ws.Range("H7").Formula = '=H5+I6'
$ws.Range("I8").Value = 9.652000000000001

$ws.Range("G9").Value = "J4"
$ws.Range("H9").Formula = '=H7+I8'
$ws.Range("I10").Value = 16.510000000000005

$ws.Range("G11").Value = "J5"
$ws.Range("H11").Formula = '=H9+I10'
$ws.Range("I12").Value = 10.794999999999987

$ws.Range("G13").Value = "J6"
$ws.Range("H13").Formula = '=H11+I12'
$ws.Range("I14").Value = 10.794999999999987

$ws.Range("G15").Value = "J7"
$ws.Range("H15").Formula = '=H13+I14'
$ws.Range("I16").Value = 9.3980000000000103

$ws.Range("G17").Value = "D2"
$ws.Range("H17").Formula = '=H15+I16'
$ws.Range("I18").Value = 9.652000000000001

$ws.Range("G19").Value = "J8"
$ws.Range("H19").Formula = '=H17+I18'

# --- columns J/K/L : rectangle centers (offset to a local 0..~120 frame) -
$ws.Range("K1").Value = 6
$ws.Range("L1").Value = 22.228000000000002
$ws.Range("J1").Formula = '=L1+K1/2'

$ws.Range("J3").Formula = '=J1+I2'
$ws.Range("K3").Value = 6
$ws.Range("L3").Formula = '=J3-K3/2'

$ws.Range("J5").Formula = '=J3+I4'
$ws.Range("K5").Value = 6
$ws.Range("L5").Formula = '=J5-K5/2'

$ws.Range("J7").Formula = '=J5+I6'
$ws.Range("K7").Value = 2.75
$ws.Range("L7").Formula = '=J7-K7/2'

$ws.Range("J9").Formula = '=J7+I8'
$ws.Range("K9").Value = 6
$ws.Range("L9").Formula = '=J9-K9/2'

$ws.Range("J11").Formula = '=J9+I10'
$ws.Range("K11").Value = 6
$ws.Range("L11").Formula = '=J11-K11/2'

$ws.Range("J13").Formula = '=J11+I12'
$ws.Range("K13").Value = 6
$ws.Range("L13").Formula = '=J13-K13/2'

$ws.Range("J15").Formula = '=J13+I14'
$ws.Range("K15").Value = 6
$ws.Range("L15").Formula = '=J15-K15/2'

$ws.Range("J17").Formula = '=J15+I16'
$ws.Range("K17").Value = 2.75
$ws.Range("L17").Formula = '=J17-K17/2'

$ws.Range("J19").Formula = '=J17+I18'
$ws.Range("K19").Value = 6
$ws.Range("L19").Formula = '=J19-K19/2'

# --- extra reference values used by the "rects" drawing helper -----------
$ws.Range("G21").Value = "rects"
$ws.Range("H21").Value = 39.918999999999997
$ws.Range("H23").Value = 97.069000000000003

# --- leave the selection where the author left it -------------------------
$ws.Range("I21").Select() | Out-Null
